# Update "想去人数" (F column) counts on the sheets that mirror the
# exhibition list ("展览" and "全部类型"). Both sheets share identical data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 155
    3  = 1775
    4  = 1674
    5  = 396
    10 = 242
    12 = 80
    16 = 32
    18 = 68
    19 = 173
    21 = 430
    22 = 324
    23 = 127
    25 = 26
    27 = 680
    28 = 2479
    31 = 501
    32 = 714
    34 = 437
    35 = 246
    36 = 374
    38 = 571
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
